$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transactions")

# Drop the custodian-based "source" naming (State Street / BNP Paribas) from the
# sample txn_id values and replace it with country/region-based naming, matching
# the simplified transaction config setup (no more 'source' assignment baked
# into the txn id strings).
$ws.Range("A6").Value = "txnid_USCust0001"
$ws.Range("A7").Value = "txnid_USCust0002"
$ws.Range("A8").Value = "txnid_USCust0003"
$ws.Range("A9").Value = "txnid_FRCust0001"
$ws.Range("A10").Value = "txnid_FRCust0002"
$ws.Range("A11").Value = "txnid_FRCust0003"

# Restore the active selection to A13 on the transactions sheet (matches the
# author's saved view state after making the edit).
$ws.Activate()
$ws.Range("A13").Select()
